# Update for MathNet.Numerics library
# Adds a new row (row 10) to Sheet1 documenting the MathNet.Numerics
# component's licence info, matching the existing "Toolkit / Component /
# Licence / Link / Comment" table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
# Column A = Toolkit, C = Licence, D = Link (B/E left blank, same as the
# sparse rows already present in the sheet, e.g. row 9).
$ws.Range("A10").Value = "MathNet.Numerics"
$ws.Range("C10").Value = "MIT/X11"
$ws.Range("D10").Value = "https://numerics.mathdotnet.com/License.html"

# --- Formatting ---------------------------------------------------------
# Give the new Toolkit/Licence cells a thin left+right border (matches the
# vertical rules used throughout the rest of the table). Apply per-cell so
# each one gets its own complete left+right rule, not just the outer edges
# of a combined range.
foreach ($addr in "A10", "C10") {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# --- Selection ------------------------------------------------------------
# Leave the cursor where the author last left it when saving.
[void]$ws.Range("E16").Select()

Write-Output "Added MathNet.Numerics row"
